# Add new payment row 09876543 (Cash) 2025-08-18T18:04:18
# and normalize the previous "latest" row's phone number to a plain number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to be the most-recent payment (row 54) had its phone
# number stored as text "09876543" to preserve the leading zero. Now that a
# newer payment is being appended, that old row's phone number collapses to
# a plain numeric value (losing the leading zero).
$ws.Cells.Item(54, 1).Value = 9876543

# Append the new payment as row 55.
# Column A keeps the leading zero, so it must be stored as text; a leading
# apostrophe forces Excel to treat the numeric-looking string as text.
$ws.Cells.Item(55, 1).Value = "'09876543"
$ws.Cells.Item(55, 3).Value = "Cash"
$ws.Cells.Item(55, 4).Value = "2025-08-18T18:04:18"
$ws.Cells.Item(55, 5).Value = 120
$ws.Cells.Item(55, 7).Value = 120
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 0
